# Atualizei dados bibi e add
# Update the 2025Q3 row (row 29) metrics in the recorrencia trimestral sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C29").Value = 119
$ws.Range("D29").Value = 22
$ws.Range("E29").Value = 97
$ws.Range("F29").Value = 3.78657487091222
